$wb = $excel.ActiveWorkbook

# --- demography sheet (sheet1 / rId1) ---
$demography = $wb.Worksheets.Item("demography")

# Correct survival value for scenario 0
$demography.Range("E2").Value = 0.7

# New "Description" column explaining each growth scenario
$demography.Range("J1").Value = "Description"
$demography.Range("J2").Value = "Null"
$demography.Range("J3").Value = "Stable"
$demography.Range("J4").Value = "Decreasing"
$demography.Range("J5").Value = "Increasing"

# --- sampling sheet (sheet3 / rId3) ---
$sampling = $wb.Worksheets.Item("sampling")

# New "Description" column explaining each sampling-effort scenario
$sampling.Range("H1").Value = "Description"
$sampling.Range("H3").Value = "Reality + 100% effort 2025-2026"
$sampling.Range("H4").Value = "Reality + 100% effort 2025-2027"
$sampling.Range("H5").Value = "Reality + 100% effort 2025-2028"
$sampling.Range("H6").Value = "Reality + 75% effort 2025-2026"
$sampling.Range("H2").Value = "100% effort 2023-2028"
$sampling.Range("H9").Value = "100% effort 2023-2025"
$sampling.Range("H8").Value = "Reality + 75% effort through 2028"
$sampling.Range("H7").Value = "Reality + 75% effort through 2027"

# New "Total Effort" column summing yearly effort per scenario
$sampling.Range("I1").Value = "Total Effort"
$sampling.Range("I2").Formula = "=SUM(B2:G2)"
$sampling.Range("I3:I9").Formula = "=SUM(B3:G3)"

# Widen the new Description column to fit its contents
$sampling.Columns.Item(8).AutoFit()

# Selection / active-sheet bookkeeping to match the saved UI state
$demography.Range("J6").Select()
$sampling.Activate()
$sampling.Range("H2:H9").Select()
